$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new data row at row 212 -------------------------------
# (shifts the existing row 212 and everything below it down by one row)
$ws.Rows(212).Insert()

$ws.Range("A212").Value  = 4
$ws.Range("B212").Value  = "Feria Lagunitas de Puerto Montt"
$ws.Range("C212").Value  = "Los Lagos"
$ws.Range("D212").Value2 = 45006
$ws.Range("E212").Value  = 10
$ws.Range("F212").Value  = 100112009
$ws.Range("G212").Value  = "Acelga"
$ws.Range("H212").Value  = "Sin especificar"
$ws.Range("I212").Value  = "Primera"
$ws.Range("J212").Value  = 80
$ws.Range("K212").Value  = 10000
$ws.Range("L212").Value  = 10000
$ws.Range("M212").Value  = 10000
$ws.Range("N212").Value  = "$/docena de atados (12 kilos)"
$ws.Range("O212").Value  = "Región de La Araucanía"
$ws.Range("P212").Value  = 833
$ws.Range("Q212").Value  = 12
$ws.Range("R212").Value  = "Hortaliza"

# --- Insert second new data row at row 242 ------------------------------
# (row numbers below already account for the insert above)
$ws.Rows(242).Insert()

$ws.Range("A242").Value  = 4
$ws.Range("B242").Value  = "Feria Lagunitas de Puerto Montt"
$ws.Range("C242").Value  = "Los Lagos"
$ws.Range("D242").Value2 = 45005
$ws.Range("E242").Value  = 10
$ws.Range("F242").Value  = 100112009
$ws.Range("G242").Value  = "Acelga"
$ws.Range("H242").Value  = "Sin especificar"
$ws.Range("I242").Value  = "Primera"
$ws.Range("J242").Value  = 25
$ws.Range("K242").Value  = 10000
$ws.Range("L242").Value  = 10000
$ws.Range("M242").Value  = 10000
$ws.Range("N242").Value  = "$/docena de atados (12 kilos)"
$ws.Range("O242").Value  = "Región de La Araucanía"
$ws.Range("P242").Value  = 833
$ws.Range("Q242").Value  = 12
$ws.Range("R242").Value  = "Hortaliza"

Write-Output ("Final used range: " + $ws.UsedRange.Address())
